$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Consumption"
$ws.Range("B3").Value = "Private Consumption Expenditure"
$ws.Range("C3").Value = "Activity"

$ws.Range("B4").Value = "Gross Capital Formation"
$ws.Range("A4").Value = "Investment"
$ws.Range("C4").Value = "Activity"

$ws.Range("E7").Select()
